$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" (F1) and "is_enabled_lbl" (G1) header columns.
# Their neighbours "order_by" (H1) and "rem" (I1) shift left into F1/G1.
$orderByText = $ws.Range("H1").Text
$remText = $ws.Range("I1").Text

$ws.Range("F1").Value = $orderByText
$ws.Range("G1").Value = $remText
$ws.Range("H1").Value = $null
$ws.Range("I1").Value = $null
